# Reworked BaseThermo and Nasa: the "geometry" column (with 'linear'/
# 'nonlinear' values) is removed from the sheet, and the "thermo_model"
# header is renamed to "statmech_model".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the entire "geometry" column (column H) - this shifts every
# column to its right one place to the left.
$ws.Range("H:H").Delete()

# Rename the thermo_model header (now still in column F) to statmech_model.
$ws.Range("F1").Value = "statmech_model"

# Match the final selection left by the edit.
$ws.Range("K8").Select()
